# Renamed few transcripts. Updated the DataSheet
# Abbreviate the "Speaker" column (D) entries throughout the transcript:
#   TRACY LEWIS -> T
#   STUDENT     -> S   (row 45 is a special case -> SN)
#   STUDENTS    -> SS
# Also update the "Student Tag" text in G45 to match the new abbreviation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 147

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)   # column D = Speaker
    $val = $cell.Text

    if ($r -eq 45) {
        $cell.Value = "SN"
    } elseif ($val -eq "TRACY LEWIS") {
        $cell.Value = "T"
    } elseif ($val -eq "STUDENTS") {
        $cell.Value = "SS"
    } elseif ($val -eq "STUDENT") {
        $cell.Value = "S"
    }
}

# Row 45's Student Tag (column G) referenced the full word "Student" - trim it too.
$ws.Cells.Item(45, 7).Value = "2 - Relating to Another S"
